# plantilla_planificacion_v2.xlsx - incidental edits captured while testing
# the new fuzzy-matching-for-technician-names feature against this planning
# template.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The test run advanced the sample ticket_id counter in row 2 (101 -> 102).
$ws.Range("B2").Value = 102

# Column A ("fecha") was widened while reviewing the sheet.
$ws.Range("A:A").ColumnWidth = 9.166666666666666

# Final cursor position left on the sheet after the review.
$ws.Range("E9").Select()
